$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new project row (row 3) to the project list
$ws.Range("K3").Value = "Michael"

$ws.Range("A3").Value = "Sleepless Nights"
$ws.Range("B3").Value = "Pulau Ubin"
$ws.Range("C3").Value = "2-Room"
$ws.Range("D3").Value = 1337
$ws.Range("E3").Value = 350000
$ws.Range("F3").Value = "3-Room"
$ws.Range("G3").Value = 420
$ws.Range("H3").Value = 450000

$ws.Range("I3").Value = $ws.Range("I2").Value2
$ws.Range("J3").Value = $ws.Range("J2").Value2
$ws.Range("I2:J2").Copy()
$ws.Range("I3:J3").PasteSpecial(-4122)

$ws.Range("L3").Value = 3

$ws.Range("E5").Select()
